# Update 2.1.1.1e workbook: add 2020 data column (N) and revise several
# existing 2019 (M) / 2018 (L) figures to match newly published data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revise existing M (2019) column values ---------------------------------
$ws.Range("M5").Value  = 68.400000000000006
$ws.Range("M6").Value  = 108.2
$ws.Range("M7").Value  = 51.7
$ws.Range("M8").Value  = 97.7
$ws.Range("M9").Value  = 106.7
$ws.Range("M10").Value = 124.2
$ws.Range("M11").Value = 138.80000000000001
$ws.Range("M12").Value = 33.9
$ws.Range("M13").Value = 96
$ws.Range("M14").Value = 7.7

# --- Revise existing L (2018) column values ---------------------------------
$ws.Range("L9").Value  = 105.6
$ws.Range("L12").Value = 27.1

# --- Seed the new N column's formatting from the current last column (M) ---
# so it inherits the same style pattern already used across the table.
$ws.Range("M3:M14").Copy()
$ws.Range("N3:N14").PasteSpecial(-4122)

# Row 3 (the separator row above the header) keeps only a bottom border with
# the plain data font/no number-format -- matching the look of the other
# border-only cells in the sheet (e.g. A14:C14) rather than M3's style.
$ws.Range("C14").Copy()
$ws.Range("N3").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# --- Populate the new N (2020) column ---------------------------------------
$ws.Range("N4").Value  = 2020
$ws.Range("N5").Value  = 68.5
$ws.Range("N6").Value  = 106.7
$ws.Range("N7").Value  = 53.2
$ws.Range("N8").Value  = 49.6
$ws.Range("N9").Value  = 108.9
$ws.Range("N10").Value = 107.8
$ws.Range("N11").Value = 155.69999999999999
$ws.Range("N12").Value = 25.9
$ws.Range("N13").Value = 103.5
$ws.Range("N14").Value = 11

# --- Page setup --------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
